# Apply the latest cryptos-list price/volume refresh (GitHub Actions data pull).
# For numeric-looking text in column D we keep it as text (leading "'" = quote-prefix,
# same as the source workbook's inlineStr cells) so Excel does not coerce it to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.440.83"
$ws.Range("E2").Value = "  -1.59%  "

$ws.Range("D3").Value = "1.840.37"
$ws.Range("E3").Value = "  -2.07%  "

$ws.Range("D4").Value = "'1.000"

$ws.Range("D5").Value = "'260.29"
$ws.Range("E5").Value = "  -6.50%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").Value = "'0.5202"
$ws.Range("E7").Value = "  -2.27%  "

$ws.Range("D8").Value = "'0.3269"
$ws.Range("E8").Value = "  -5.03%  "

$ws.Range("D9").Value = "'0.06779"
$ws.Range("E9").Value = "  -2.62%  "

$ws.Range("D10").Value = "'18.59"
$ws.Range("E10").Value = "  -7.34%  "

$ws.Range("D11").Value = "'0.7665"
$ws.Range("E11").Value = "  -4.79%  "

$ws.Range("D12").Value = "'0.07699"
$ws.Range("E12").Value = "  -0.63%  "

$ws.Range("D13").Value = "1.820.00"
$ws.Range("E13").Value = "  -3.27%  "

$ws.Range("D14").Value = "'88.11"
$ws.Range("E14").Value = "  -2.61%  "

$ws.Range("D15").Value = "'5.030"
$ws.Range("E15").Value = "  -2.73%  "

$ws.Range("D16").Value = "'0.9999"
$ws.Range("E16").Value = "  -0.06%  "

$ws.Range("D17").Value = "'13.90"
$ws.Range("E17").Value = "  -4.50%  "

$ws.Range("D18").Value = "'1.000"
$ws.Range("E18").Value = "  -0.05%  "

$ws.Range("E19").Value = "  -1.06%  "

$ws.Range("D20").Value = "26.446.74"
$ws.Range("E20").Value = "  -1.73%  "

$ws.Range("D21").Value = "2.070.46"
$ws.Range("E21").Value = "  -2.36%  "

$ws.Range("E22").Value = "  -4.03%  "

$ws.Range("D23").Value = "'9.475"
$ws.Range("E23").Value = "  -5.60%  "

$ws.Range("D24").Value = "'5.966"
$ws.Range("E24").Value = "  -3.88%  "

$ws.Range("D25").Value = "'144.39"
$ws.Range("E25").Value = "  -1.72%  "

$ws.Range("D26").Value = "'2.219"
$ws.Range("E26").Value = "  -6.40%  "

$ws.Range("E27").Value = "  -0.94%  "

$ws.Range("D28").Value = "'16.99"
$ws.Range("E28").Value = "  -2.08%  "

$ws.Range("D29").Value = "'111.10"
$ws.Range("E29").Value = "  -2.31%  "

$ws.Range("D30").Value = "'4.173"
$ws.Range("E30").Value = "  -4.16%  "

$ws.Range("D31").Value = "'4.126"
$ws.Range("E31").Value = "  -4.34%  "

$ws.Range("D32").Value = "'0.08713"
$ws.Range("E32").Value = "  -2.08%  "

$ws.Range("D33").Value = "'0.04797"
$ws.Range("E33").Value = "  -2.12%  "

$ws.Range("D34").Value = "'1.126"
$ws.Range("E34").Value = "  -4.30%  "

$ws.Range("D35").Value = "'2.836"
$ws.Range("E35").Value = "  -1.89%  "

$ws.Range("D36").Value = "'0.7007"
$ws.Range("E36").Value = "  -3.51%  "

$ws.Range("D37").Value = "'3.073"
$ws.Range("E37").Value = "  -6.42%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'2.208"
$ws.Range("E38").Value = "  -6.46%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01752"
$ws.Range("E39").Value = "  -5.15%  "

$ws.Range("D40").Value = "'0.4824"
$ws.Range("E40").Value = "  -5.90%  "

$ws.Range("D41").Value = "'110.91"
$ws.Range("E41").Value = "  -4.33%  "

$ws.Range("D42").Value = "'0.8887"
$ws.Range("E42").Value = "  -7.14%  "

$ws.Range("D43").Value = "'6.067"
$ws.Range("E43").Value = "  -1.91%  "

$ws.Range("D44").Value = "'0.9999"
$ws.Range("E44").Value = "  -0.05%  "

$ws.Range("D45").Value = "'7.684"
$ws.Range("E45").Value = "  -5.22%  "

$ws.Range("D46").Value = "'0.05870"
$ws.Range("E46").Value = "  -1.36%  "

$ws.Range("D47").Value = "'0.4119"
$ws.Range("E47").Value = "  -7.81%  "

$ws.Range("D48").Value = "'8.941"
$ws.Range("E48").Value = "  -4.50%  "

$ws.Range("E49").Value = "  -3.42%  "

$ws.Range("D50").Value = "'0.1218"
$ws.Range("E50").Value = "  -9.19%  "

$ws.Range("D51").Value = "'0.8854"
$ws.Range("E51").Value = "  +0.05%  "
